$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values for the symbol list refresh.
# Values must be written as Text (matching the workbook's existing text-stored
# numbers/percentages), so we force NumberFormat = "@" before assigning, then
# reset the style back to Normal so no stray number-format style is left behind.
$cellUpdates = @{
    "D2" = "299.86"
    "E2" = "-2.50%"
    "D3" = "40.71"
    "E3" = "-0.71%"
    "D4" = "5.143"
    "E4" = "-1.28%"
    "D5" = "0.07509"
    "E5" = "-2.08%"
    "E6" = "0.92%"
    "D7" = "1.629"
    "E7" = "-0.18%"
    "D8" = "0.9392"
    "E8" = "2.52%"
    "D10" = "0.1198"
    "E10" = "-1.13%"
    "D11" = "0.1791"
    "E11" = "-1.26%"
    "D12" = "0.08925"
    "E12" = "-2.02%"
    "D13" = "0.04150"
    "E13" = "-0.06%"
    "D14" = "0.1054"
    "E14" = "0.22%"
    "D15" = "0.001268"
    "E15" = "0.79%"
    "D16" = "0.005793"
    "E16" = "1.10%"
    "D17" = "3.346"
    "D18" = "0.3353"
    "E18" = "0.52%"
    "D19" = "7.653"
    "E19" = "3.73%"
    "E20" = "-1.76%"
    "D21" = "0.2815"
    "E21" = "3.73%"
    "D22" = "0.03840"
    "E22" = "-4.42%"
    "D23" = "0.001286"
    "E23" = "2.19%"
    "D24" = "0.003960"
    "E24" = "-7.06%"
    "D25" = "0.0001305"
    "E25" = "0.40%"
    "D38" = "0.02347"
    "E38" = "-5.68%"
    "D39" = "0.05079"
    "E39" = "-4.92%"
    "D40" = "0.007772"
    "E40" = "-0.99%"
    "D41" = "0.1295"
    "E41" = "-1.46%"
    "D42" = "0.007583"
    "E42" = "16.54%"
    "D43" = "0.003567"
    "E43" = "86.43%"
    "D44" = "0.008003"
    "E44" = "-3.06%"
    "D45" = "0.3252"
    "E45" = "-2.34%"
    "D46" = "0.00006823"
    "E46" = "1.39%"
    "E47" = "0.22%"
    "D48" = "0.2516"
    "E48" = "-6.85%"
    "E49" = "35.84%"
    "D50" = "0.00002108"
    "E50" = "0.22%"
    "D51" = "0.0002008"
    "E51" = "0.22%"
}

foreach ($ref in $cellUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$ref]
    $cell.Style = "Normal"
}
